$wb = $excel.ActiveWorkbook

# Rename sheets: "mars" -> "March", "april" -> "April"
$wsMarch = $wb.Worksheets.Item("mars")
$wsMarch.Name = "March"

$wsApril = $wb.Worksheets.Item("april")
$wsApril.Name = "April"

# Add two new transaction rows to the April sheet and update the monthly total
$wsApril.Range("A7").Value = "Transportation"
$wsApril.Range("B7").Value = "rdfr"
$wsApril.Range("C7").Value = "'2023-04-12"
$wsApril.Range("C7").ClearFormats()
$wsApril.Range("D7").Value = 656.0
$wsApril.Range("E7").Value = "Checkings"

$wsApril.Range("A8").Value = "Other"
$wsApril.Range("B8").Value = "emsd"
$wsApril.Range("C8").Value = "'2023-04-12"
$wsApril.Range("C8").ClearFormats()
$wsApril.Range("D8").Value = 450.0
$wsApril.Range("E8").Value = "Savings"

$wsApril.Range("A9").Value = "Monthly total: "
$wsApril.Range("B9").Value = 3212.0
